$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells used to clone number-format/style onto newly created cells,
# since newly-materialized cells in this engine inherit the rows previously-seen
# style rather than the column default style.
$refA = $ws.Range("A12")
$refB = $ws.Range("B9")
$refC = $ws.Range("C9")

# --- Apply new/changed cell values ---
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "144651 - Antonio Fernando Sartori"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "144651 - Antonio Fernando Sartori"

$refA.Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "Programa resumido:"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Semestral"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Short syllabus:"

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "Programa:"

$refB.Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"

$refC.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Syllabus:"

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "Método:"

$refB.Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "144651 - Antonio Fernando Sartori"

$refC.Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "144651 - Antonio Fernando Sartori"

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Bibliografia:"

$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "Requisitos:"

$refB.Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)`n"

$refC.Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)`n"

# --- Clear cells whose content was removed ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Delete the now-removed trailing row (old row 24) ---
$ws.Rows.Item(24).EntireRow.Delete()

# --- Fix up row heights to match the target layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(22).RowHeight = 15
